$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1229006.6
$ws.Range("J17").Value = 1263124.9
$ws.Range("L17").Value = 3789374.7
$ws.Range("N17").Value = -3789710.7

$ws.Range("H19").Value = 1171.3158
$ws.Range("I19").Value = 677
$ws.Range("K19").Value = 677
$ws.Range("M19").Value = -502

$ws.Range("H28").Value = 4238.7856
$ws.Range("J28").Value = 2548
$ws.Range("L28").Value = 2548
$ws.Range("N28").Value = -3518

$ws.Range("H32").Value = 6169.6
$ws.Range("J32").Value = 5116
$ws.Range("L32").Value = 5116
$ws.Range("N32").Value = -5768

$ws.Range("H33").Value = 576.7857
$ws.Range("I33").Value = 290.91666
$ws.Range("J33").Value = 2292
$ws.Range("K33").Value = 290.91666
$ws.Range("L33").Value = 2292
$ws.Range("M33").Value = -61.91665999999998
$ws.Range("N33").Value = -2750

$ws.Range("H74").Value = 4997.3477
$ws.Range("I74").Value = 3996
$ws.Range("K74").Value = 3996
$ws.Range("M74").Value = -3060

$ws.Range("H77").Value = 4997.3477
$ws.Range("I77").Value = 3996
$ws.Range("K77").Value = 19980
$ws.Range("M77").Value = -15300

$ws.Range("H116").Value = 1858218.9
$ws.Range("I116").Value = 3707705.2
$ws.Range("K116").Value = 3707705.2
$ws.Range("M116").Value = -3704263.2

$ws.Range("H132").Value = 3585.6829
$ws.Range("I132").Value = 3861.5454
$ws.Range("J132").Value = 2447.75
$ws.Range("K132").Value = 11584.6362
$ws.Range("L132").Value = 7343.25
$ws.Range("M132").Value = -9054.636200000001
$ws.Range("N132").Value = -12403.25

$ws.Range("H137").Value = 1187472.5
$ws.Range("I137").Value = 1511750.6
$ws.Range("K137").Value = 4535251.800000001
$ws.Range("M137").Value = -4532701.800000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H9").Value = 20008
$ws.Range("I9").Value = 20008
$ws.Range("K9").Value = 20008
$ws.Range("M9").Value = -19838

$ws.Range("H20").Value = 20008
$ws.Range("I20").Value = 20008
$ws.Range("K20").Value = 20008
$ws.Range("M20").Value = -19738

$ws.Range("H32").Value = 2706.8635
$ws.Range("I32").Value = 2645.3096
$ws.Range("J32").Value = 3999.5
$ws.Range("K32").Value = 2645.3096
$ws.Range("L32").Value = 3999.5
$ws.Range("M32").Value = -2358.3096
$ws.Range("N32").Value = -4573.5

$ws.Range("H45").Value = 157697.22
$ws.Range("I45").Value = 240667.89
$ws.Range("K45").Value = 240667.89
$ws.Range("M45").Value = -240290.89

$ws.Range("H61").Value = 6419.5
$ws.Range("I61").Value = 6769.0713
$ws.Range("K61").Value = 6769.0713
$ws.Range("M61").Value = -6557.0713

$ws.Range("H74").Value = 3864.24
$ws.Range("I74").Value = 2072.8333
$ws.Range("K74").Value = 2072.8333
$ws.Range("M74").Value = -1198.8333

$ws.Range("H77").Value = 3864.24
$ws.Range("I77").Value = 2072.8333
$ws.Range("K77").Value = 10364.1665
$ws.Range("M77").Value = -5996.166499999999

$ws.Range("H110").Value = 0
$ws.Range("I110").Value = 0
$ws.Range("K110").Value = 0
$ws.Range("M110").ClearContents()

$ws.Range("H122").Value = 6289428
$ws.Range("I122").Value = 4998
$ws.Range("J122").Value = 7336833
$ws.Range("K122").Value = 14994
$ws.Range("L122").Value = 22010499
$ws.Range("M122").Value = -12544
$ws.Range("N122").Value = -22015399

$ws.Range("H132").Value = 2822.0303
$ws.Range("I132").Value = 2234
$ws.Range("J132").Value = 4174.5
$ws.Range("K132").Value = 6702
$ws.Range("L132").Value = 12523.5
$ws.Range("M132").Value = -4172
$ws.Range("N132").Value = -17583.5

$ws.Range("H136").Value = 6419.5
$ws.Range("I136").Value = 6769.0713
$ws.Range("K136").Value = 20307.2139
$ws.Range("M136").Value = -17757.2139

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H42").Value = 204755
$ws.Range("J42").Value = 204755
$ws.Range("L42").Value = 204755
$ws.Range("N42").Value = -205411

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 688.2857
$ws.Range("I22").Value = 685.8
$ws.Range("J22").Value = 689.6667
$ws.Range("K22").Value = 685.8
$ws.Range("L22").Value = 689.6667
$ws.Range("M22").Value = -335.8
$ws.Range("N22").Value = -1389.6667

$ws.Range("H58").Value = 3522.9048
$ws.Range("I58").Value = 2199
$ws.Range("K58").Value = 2199
$ws.Range("M58").Value = -1996

$ws.Range("H105").Value = 106742.7
$ws.Range("I105").Value = 151532.64
$ws.Range("J105").Value = 2232.8333
$ws.Range("K105").Value = 151532.64
$ws.Range("L105").Value = 2232.8333
$ws.Range("M105").Value = -149785.64
$ws.Range("N105").Value = -5726.8333

$ws.Range("H107").Value = 125042504
$ws.Range("I107").Value = 200066610
$ws.Range("J107").Value = 2333
$ws.Range("K107").Value = 200066610
$ws.Range("L107").Value = 2333
$ws.Range("M107").Value = -200064690
$ws.Range("N107").Value = -6173

$ws.Range("H132").Value = 16516.027
$ws.Range("I132").Value = 5579.9355
$ws.Range("K132").Value = 16739.8065
$ws.Range("M132").Value = -14209.8065

$ws.Range("H136").Value = 3522.9048
$ws.Range("I136").Value = 2199
$ws.Range("K136").Value = 6597
$ws.Range("M136").Value = -4047

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H22").Value = 500
$ws.Range("I22").Value = 500
$ws.Range("K22").Value = 1500
$ws.Range("M22").Value = -1331

$ws.Range("H27").Value = 500
$ws.Range("I27").Value = 500
$ws.Range("K27").Value = 1500
$ws.Range("M27").Value = -1398

$ws.Range("H37").Value = 59162.883
$ws.Range("J37").Value = 59162.883
$ws.Range("L37").Value = 177488.649
$ws.Range("N37").Value = -177712.649

$ws.Range("H97").Value = 38865.75
$ws.Range("J97").Value = 1791.8889
$ws.Range("L97").Value = 5375.6667
$ws.Range("N97").Value = -6367.6667

$ws.Range("H98").Value = 629.6
$ws.Range("I98").Value = 617.6667
$ws.Range("J98").Value = 647.5
$ws.Range("K98").Value = 1853.0001
$ws.Range("L98").Value = 1942.5
$ws.Range("M98").Value = -355.0001
$ws.Range("N98").Value = -4938.5

$ws.Range("H139").Value = 3336744.5
$ws.Range("I139").Value = 10002400
$ws.Range("J139").Value = 3916.6667
$ws.Range("K139").Value = 30007200
$ws.Range("L139").Value = 11750.0001
$ws.Range("M139").Value = -30002060
$ws.Range("N139").Value = -22030.0001

$ws.Range("H140").Value = 9245.862999999999
$ws.Range("I140").Value = 9638.714
$ws.Range("J140").Value = 996
$ws.Range("K140").Value = 28916.142
$ws.Range("L140").Value = 2988
$ws.Range("M140").Value = -23736.142
$ws.Range("N140").Value = -13348

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H3").Value = 618105.1
$ws.Range("J3").Value = 669572.25
$ws.Range("L3").Value = 669572.25
$ws.Range("N3").Value = -669804.25

$ws.Range("H102").Value = 12390.917
$ws.Range("I102").Value = 15410.111
$ws.Range("K102").Value = 15410.111
$ws.Range("M102").Value = -13788.111

$ws.Range("H132").Value = 2567.1035
$ws.Range("I132").Value = 1555.3334
$ws.Range("K132").Value = 4666.0002
$ws.Range("M132").Value = -2136.0002

$ws.Range("H141").Value = 99993.5
$ws.Range("J141").Value = 99993.5
$ws.Range("L141").Value = 99993.5
$ws.Range("N141").Value = -110353.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 27496.475
$ws.Range("I40").Value = 37083.332
$ws.Range("J40").Value = 11061.857
$ws.Range("K40").Value = 37083.332
$ws.Range("L40").Value = 11061.857
$ws.Range("M40").Value = -36947.332
$ws.Range("N40").Value = -11333.857

$ws.Range("H46").Value = 3538.1177
$ws.Range("I46").Value = 1151.4
$ws.Range("K46").Value = 1151.4
$ws.Range("M46").Value = -963.4000000000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H103").Value = 0
$ws.Range("J103").Value = 0
$ws.Range("L103").Value = 0
$ws.Range("N103").ClearContents()

$ws.Range("H107").Value = 19727.824
$ws.Range("I107").Value = 2267
$ws.Range("J107").Value = 44671.855
$ws.Range("K107").Value = 6801
$ws.Range("L107").Value = 134015.565
$ws.Range("M107").Value = -4881
$ws.Range("N107").Value = -137855.565

$ws.Range("H132").Value = 26499.625
$ws.Range("I132").Value = 28143.143
$ws.Range("K132").Value = 84429.429
$ws.Range("M132").Value = -81899.429

$ws.Range("H136").Value = 3226.6086
$ws.Range("I136").Value = 2234.2778
$ws.Range("K136").Value = 6702.8334
$ws.Range("M136").Value = -4152.8334
